# Update model comparison table (rows 3-6) with refreshed analysis results.
# B=Loglike, C=X2, D=G2, E=DF (unchanged), F=AIC, G=BIC

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = -2328.346493560788
$ws.Range("C3").Value = 4024.079716055674
$ws.Range("D3").Value = 480.3115075978895
$ws.Range("F3").Value = 4718.692987121575
$ws.Range("G3").Value = 4846.962525798249

$ws.Range("B4").Value = -2209.92758755932
$ws.Range("C4").Value = 3199.50935591981
$ws.Range("D4").Value = 408.7700485178675
$ws.Range("F4").Value = 4513.855175118641
$ws.Range("G4").Value = 4708.328346660694

$ws.Range("B5").Value = -2198.768944509708
$ws.Range("C5").Value = 1310.539557061597
$ws.Range("D5").Value = 423.0238552106839
$ws.Range("F5").Value = 4523.537889019416
$ws.Range("G5").Value = 4784.214693426849

$ws.Range("B6").Value = -2102.035380288189
$ws.Range("C6").Value = 1968.34373054423
$ws.Range("D6").Value = 336.241794453155
$ws.Range("F6").Value = 4362.070760576378
$ws.Range("G6").Value = 4688.951197849191
